$wb = $excel.ActiveWorkbook

# The "Gaz" sheet holds the Contract/Last/High/Low table that needs reshaping
# into a Bid/Ask/Last table (and losing the trailing column E).
$ws = $wb.Worksheets.Item("Gaz")

# Header row: Contract -> Bid, Last -> Ask, High -> Last, Low (col E) removed.
$ws.Range("B1").Value = "Bid"
$ws.Range("C1").Value = "Ask"
$ws.Range("D1").Value = "Last"

# Data row: contract name cell becomes a placeholder dash; D2 keeps the
# existing "-" that used to live in E2, and column E is cleared entirely.
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"

# Drop the now-unused column E so the used range shrinks back to A1:D2.
$ws.Columns.Item(5).Delete()
